# Auto-generated: apply year-to-date crime count updates for 2023-11-20
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = 39   # B2: 38 -> 39
$ws.Cells.Item(2, 6).Value = 85   # F2: 84 -> 85
$ws.Cells.Item(2, 10).Value = 114   # J2: 113 -> 114
$ws.Cells.Item(3, 6).Value = 128   # F3: 127 -> 128
$ws.Cells.Item(3, 8).Value = 149   # H3: 147 -> 149
$ws.Cells.Item(3, 10).Value = 216   # J3: 215 -> 216
$ws.Cells.Item(9, 2).Value = 361   # B9: 360 -> 361
$ws.Cells.Item(9, 4).Value = 388   # D9: 387 -> 388
$ws.Cells.Item(9, 5).Value = 444   # E9: 443 -> 444
$ws.Cells.Item(9, 6).Value = 495   # F9: 493 -> 495
$ws.Cells.Item(9, 7).Value = 422   # G9: 421 -> 422
$ws.Cells.Item(9, 9).Value = 478   # I9: 477 -> 478
$ws.Cells.Item(9, 10).Value = 398   # J9: 396 -> 398
$ws.Cells.Item(10, 2).Value = 1276   # B10: 1269 -> 1276
$ws.Cells.Item(10, 3).Value = 1509   # C10: 1506 -> 1509
$ws.Cells.Item(10, 4).Value = 1714   # D10: 1711 -> 1714
$ws.Cells.Item(10, 5).Value = 2056   # E10: 2050 -> 2056
$ws.Cells.Item(10, 6).Value = 2010   # F10: 2004 -> 2010
$ws.Cells.Item(10, 7).Value = 871   # G10: 870 -> 871
$ws.Cells.Item(10, 9).Value = 813   # I10: 811 -> 813
$ws.Cells.Item(10, 10).Value = 682   # J10: 681 -> 682
$ws.Cells.Item(11, 2).Value = 1765   # B11: 1756 -> 1765
$ws.Cells.Item(11, 3).Value = 2121   # C11: 2118 -> 2121
$ws.Cells.Item(11, 4).Value = 2326   # D11: 2322 -> 2326
$ws.Cells.Item(11, 5).Value = 2719   # E11: 2712 -> 2719
$ws.Cells.Item(11, 6).Value = 2729   # F11: 2719 -> 2729
$ws.Cells.Item(11, 7).Value = 1519   # G11: 1517 -> 1519
$ws.Cells.Item(11, 8).Value = 1268   # H11: 1266 -> 1268
$ws.Cells.Item(11, 9).Value = 1618   # I11: 1615 -> 1618
$ws.Cells.Item(11, 10).Value = 1440   # J11: 1435 -> 1440

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 4).Value = 18   # D5: 16 -> 18
$ws.Cells.Item(7, 2).Value = 49   # B7: 48 -> 49
$ws.Cells.Item(8, 6).Value = 139   # F8: 137 -> 139
$ws.Cells.Item(16, 5).Value = 7   # E16: 6 -> 7
$ws.Cells.Item(20, 6).Value = 13   # F20: 12 -> 13
$ws.Cells.Item(21, 2).Value = 16   # B21: 15 -> 16
$ws.Cells.Item(22, 8).Value = 4   # H22: 3 -> 4
$ws.Cells.Item(27, 3).Value = 25   # C27: 24 -> 25
$ws.Cells.Item(28, 2).Value = 94   # B28: 93 -> 94
$ws.Cells.Item(28, 6).Value = 119   # F28: 116 -> 119
$ws.Cells.Item(28, 9).Value = 84   # I28: 83 -> 84
$ws.Cells.Item(29, 6).Value = 26   # F29: 25 -> 26
$ws.Cells.Item(32, 2).Value = 60   # B32: 58 -> 60
$ws.Cells.Item(32, 10).Value = 72   # J32: 71 -> 72
$ws.Cells.Item(41, 5).Value = 26   # E41: 25 -> 26
$ws.Cells.Item(45, 2).Value = 28   # B45: 27 -> 28
$ws.Cells.Item(49, 3).Value = 13   # C49: 12 -> 13
$ws.Cells.Item(50, 10).Value = 32   # J50: 31 -> 32
$ws.Cells.Item(52, 10).Value = 27   # J52: 26 -> 27
$ws.Cells.Item(53, 2).Value = 251   # B53: 250 -> 251
$ws.Cells.Item(53, 3).Value = 358   # C53: 357 -> 358
$ws.Cells.Item(53, 5).Value = 687   # E53: 684 -> 687
$ws.Cells.Item(53, 6).Value = 597   # F53: 595 -> 597
$ws.Cells.Item(53, 9).Value = 307   # I53: 306 -> 307
$ws.Cells.Item(54, 7).Value = 10   # G54: 9 -> 10
$ws.Cells.Item(59, 4).Value = 6   # D59: 5 -> 6
$ws.Cells.Item(61, 10).Value = 4   # J61: 3 -> 4
$ws.Cells.Item(62, 9).Value = 25   # I62: 24 -> 25
$ws.Cells.Item(63, 2).Value = 14   # B63: 13 -> 14
$ws.Cells.Item(65, 5).Value = 48   # E65: 47 -> 48
$ws.Cells.Item(67, 8).Value = 9   # H67: 8 -> 9
$ws.Cells.Item(70, 10).Value = 25   # J70: 24 -> 25
$ws.Cells.Item(80, 6).Value = 35   # F80: 34 -> 35
$ws.Cells.Item(89, 4).Value = 20   # D89: 19 -> 20
$ws.Cells.Item(95, 5).Value = 99   # E95: 98 -> 99
$ws.Cells.Item(95, 7).Value = 15   # G95: 14 -> 15
$ws.Cells.Item(96, 2).Value = 17   # B96: 16 -> 17
$ws.Cells.Item(99, 2).Value = 1765   # B99: 1756 -> 1765
$ws.Cells.Item(99, 3).Value = 2121   # C99: 2118 -> 2121
$ws.Cells.Item(99, 4).Value = 2326   # D99: 2322 -> 2326
$ws.Cells.Item(99, 5).Value = 2719   # E99: 2712 -> 2719
$ws.Cells.Item(99, 6).Value = 2729   # F99: 2719 -> 2729
$ws.Cells.Item(99, 7).Value = 1519   # G99: 1517 -> 1519
$ws.Cells.Item(99, 8).Value = 1268   # H99: 1266 -> 1268
$ws.Cells.Item(99, 9).Value = 1618   # I99: 1615 -> 1618
$ws.Cells.Item(99, 10).Value = 1440   # J99: 1435 -> 1440

# Sheet 4: Edgewater
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 3).Value = 19   # C5: 18 -> 19
$ws.Cells.Item(6, 3).Value = 25   # C6: 24 -> 25

# Sheet 6: Auburn Gresham
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(6, 2).Value = 31   # B6: 30 -> 31
$ws.Cells.Item(7, 2).Value = 49   # B7: 48 -> 49

# Sheet 8: Austin
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(6, 6).Value = 31   # F6: 30 -> 31
$ws.Cells.Item(7, 6).Value = 94   # F7: 93 -> 94
$ws.Cells.Item(8, 6).Value = 139   # F8: 137 -> 139

# Sheet 10: Chinatown
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(8, 2).Value = 10   # B8: 9 -> 10
$ws.Cells.Item(9, 2).Value = 16   # B9: 15 -> 16

# Sheet 11: Chicago Lawn
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(6, 6).Value = 6   # F6: 5 -> 6
$ws.Cells.Item(7, 6).Value = 13   # F7: 12 -> 13

# Sheet 12: Garfield Park
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 10).Value = 2   # J2: 1 -> 2
$ws.Cells.Item(7, 2).Value = 12   # B7: 11 -> 12
$ws.Cells.Item(8, 2).Value = 46   # B8: 45 -> 46
$ws.Cells.Item(9, 2).Value = 60   # B9: 58 -> 60
$ws.Cells.Item(9, 10).Value = 72   # J9: 71 -> 72

# Sheet 15: Loop
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(8, 2).Value = 203   # B8: 202 -> 203
$ws.Cells.Item(8, 3).Value = 306   # C8: 305 -> 306
$ws.Cells.Item(8, 5).Value = 607   # E8: 604 -> 607
$ws.Cells.Item(8, 6).Value = 524   # F8: 522 -> 524
$ws.Cells.Item(8, 9).Value = 185   # I8: 184 -> 185
$ws.Cells.Item(9, 2).Value = 251   # B9: 250 -> 251
$ws.Cells.Item(9, 3).Value = 358   # C9: 357 -> 358
$ws.Cells.Item(9, 5).Value = 687   # E9: 684 -> 687
$ws.Cells.Item(9, 6).Value = 597   # F9: 595 -> 597
$ws.Cells.Item(9, 9).Value = 307   # I9: 306 -> 307

# Sheet 16: Armour Square
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(6, 4).Value = 11   # D6: 9 -> 11
$ws.Cells.Item(7, 4).Value = 18   # D7: 16 -> 18

# Sheet 17: Old Town
$ws = $wb.Worksheets.Item(17)
$ws.Cells.Item(3, 10).Value = 6   # J3: 5 -> 6
$ws.Cells.Item(8, 10).Value = 25   # J8: 24 -> 25

# Sheet 18: Little Italy, UIC
$ws = $wb.Worksheets.Item(18)
$ws.Cells.Item(5, 10).Value = 10   # J5: 9 -> 10
$ws.Cells.Item(7, 10).Value = 32   # J7: 31 -> 32

# Sheet 19: North Lawndale
$ws = $wb.Worksheets.Item(19)
$ws.Cells.Item(7, 5).Value = 36   # E7: 35 -> 36
$ws.Cells.Item(8, 5).Value = 48   # E8: 47 -> 48

# Sheet 20: Washington Park
$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(5, 4).Value = 10   # D5: 9 -> 10
$ws.Cells.Item(6, 4).Value = 20   # D6: 19 -> 20

# Sheet 21: Sheffield &amp; DePaul
$ws = $wb.Worksheets.Item(21)
$ws.Cells.Item(6, 6).Value = 26   # F6: 25 -> 26
$ws.Cells.Item(7, 6).Value = 35   # F7: 34 -> 35

# Sheet 22: Humboldt Park
$ws = $wb.Worksheets.Item(22)
$ws.Cells.Item(4, 5).Value = 6   # E4: 5 -> 6
$ws.Cells.Item(6, 5).Value = 26   # E6: 25 -> 26

# Sheet 26: Englewood
$ws = $wb.Worksheets.Item(26)
$ws.Cells.Item(2, 6).Value = 7   # F2: 6 -> 7
$ws.Cells.Item(3, 6).Value = 10   # F3: 9 -> 10
$ws.Cells.Item(7, 9).Value = 20   # I7: 19 -> 20
$ws.Cells.Item(8, 2).Value = 61   # B8: 60 -> 61
$ws.Cells.Item(8, 6).Value = 65   # F8: 64 -> 65
$ws.Cells.Item(9, 2).Value = 94   # B9: 93 -> 94
$ws.Cells.Item(9, 6).Value = 119   # F9: 116 -> 119
$ws.Cells.Item(9, 9).Value = 84   # I9: 83 -> 84

# Sheet 28: Jefferson Park
$ws = $wb.Worksheets.Item(28)
$ws.Cells.Item(2, 2).Value = 1   # B2: None -> 1
$ws.Cells.Item(7, 2).Value = 28   # B7: 27 -> 28

# Sheet 29: Fuller Park
$ws = $wb.Worksheets.Item(29)
$ws.Cells.Item(7, 6).Value = 11   # F7: 10 -> 11
$ws.Cells.Item(9, 6).Value = 26   # F9: 25 -> 26

# Sheet 34: Norwood Park
$ws = $wb.Worksheets.Item(34)
$ws.Cells.Item(2, 8).Value = 1   # H2: None -> 1
$ws.Cells.Item(6, 8).Value = 9   # H6: 8 -> 9

# Sheet 38: Logan Square
$ws = $wb.Worksheets.Item(38)
$ws.Cells.Item(6, 10).Value = 7   # J6: 6 -> 7
$ws.Cells.Item(8, 10).Value = 27   # J8: 26 -> 27

# Sheet 40: Near South Side
$ws = $wb.Worksheets.Item(40)
$ws.Cells.Item(6, 9).Value = 12   # I6: 11 -> 12
$ws.Cells.Item(7, 9).Value = 25   # I7: 24 -> 25

# Sheet 43: Lower West Side
$ws = $wb.Worksheets.Item(43)
$ws.Cells.Item(4, 7).Value = 3   # G4: 2 -> 3
$ws.Cells.Item(6, 7).Value = 10   # G6: 9 -> 10

# Sheet 51: Morgan Park
$ws = $wb.Worksheets.Item(51)
$ws.Cells.Item(4, 4).Value = 4   # D4: 3 -> 4
$ws.Cells.Item(6, 4).Value = 6   # D6: 5 -> 6

# Sheet 52: Lincoln Square
$ws = $wb.Worksheets.Item(52)
$ws.Cells.Item(5, 3).Value = 11   # C5: 10 -> 11
$ws.Cells.Item(6, 3).Value = 13   # C6: 12 -> 13

# Sheet 54: Clearing
$ws = $wb.Worksheets.Item(54)
$ws.Cells.Item(3, 8).Value = 1   # H3: None -> 1
$ws.Cells.Item(7, 8).Value = 4   # H7: 3 -> 4

# Sheet 55: West Town
$ws = $wb.Worksheets.Item(55)
$ws.Cells.Item(6, 5).Value = 88   # E6: 87 -> 88
$ws.Cells.Item(6, 7).Value = 11   # G6: 10 -> 11
$ws.Cells.Item(7, 5).Value = 99   # E7: 98 -> 99
$ws.Cells.Item(7, 7).Value = 15   # G7: 14 -> 15

# Sheet 56: Wicker Park
$ws = $wb.Worksheets.Item(56)
$ws.Cells.Item(6, 2).Value = 17   # B6: 16 -> 17
$ws.Cells.Item(7, 2).Value = 17   # B7: 16 -> 17

# Sheet 58: New City
$ws = $wb.Worksheets.Item(58)
$ws.Cells.Item(5, 2).Value = 8   # B5: 7 -> 8
$ws.Cells.Item(6, 2).Value = 14   # B6: 13 -> 14

# Sheet 71: Bucktown
$ws = $wb.Worksheets.Item(71)
$ws.Cells.Item(5, 5).Value = 5   # E5: 4 -> 5
$ws.Cells.Item(6, 5).Value = 7   # E6: 6 -> 7
